$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4248.6665
$ws.Range("I51").Value = 3873
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 3873
$ws.Range("L51").Value = 5000
$ws.Range("M51").Value = -3389
$ws.Range("N51").Value = -5968

$ws.Range("H70").Value = 1494.4445
$ws.Range("J70").Value = 1500
$ws.Range("L70").Value = 4500
$ws.Range("N70").Value = -5040

$ws.Range("H73").Value = 1494.4445
$ws.Range("J73").Value = 1500
$ws.Range("L73").Value = 4500
$ws.Range("N73").Value = -6372

$ws.Range("H100").Value = 2770
$ws.Range("I100").Value = 2770
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2770
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -2229

$ws.Range("H107").Value = 2601.8462
$ws.Range("I107").Value = 2806.182
$ws.Range("J107").Value = 1478
$ws.Range("K107").Value = 2806.182
$ws.Range("L107").Value = 1478
$ws.Range("M107").Value = -886.1819999999998
$ws.Range("N107").Value = -5318

$ws.Range("H141").Value = 764.6087
$ws.Range("I141").Value = 753.9091
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 2261.7273
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = 2918.2727
$ws.Range("N141").Value = -13360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6475.091
$ws.Range("I2").Value = 2841.8
$ws.Range("J2").Value = 9502.833000000001
$ws.Range("K2").Value = 2841.8
$ws.Range("L2").Value = 9502.833000000001
$ws.Range("M2").Value = -2728.8
$ws.Range("N2").Value = -9728.833000000001

$ws.Range("H32").Value = 4201.815
$ws.Range("I32").Value = 2555.6191
$ws.Range("J32").Value = 9963.5
$ws.Range("K32").Value = 2555.6191
$ws.Range("L32").Value = 9963.5
$ws.Range("M32").Value = -2268.6191
$ws.Range("N32").Value = -10537.5

$ws.Range("H45").Value = 2802.8
$ws.Range("I45").Value = 1263.1428
$ws.Range("K45").Value = 1263.1428
$ws.Range("M45").Value = -886.1428000000001

$ws.Range("H61").Value = 906.3333
$ws.Range("I61").Value = 609.5
$ws.Range("K61").Value = 609.5
$ws.Range("M61").Value = -397.5

$ws.Range("H74").Value = 1153.2941
$ws.Range("I74").Value = 1153.2941
$ws.Range("K74").Value = 1153.2941
$ws.Range("M74").Value = -279.2941000000001

$ws.Range("H77").Value = 1153.2941
$ws.Range("I77").Value = 1153.2941
$ws.Range("K77").Value = 5766.4705
$ws.Range("M77").Value = -1398.4705

$ws.Range("H102").Value = 3632.6667
$ws.Range("I102").Value = 1949
$ws.Range("K102").Value = 1949
$ws.Range("M102").Value = -327

$ws.Range("H116").Value = 6475.091
$ws.Range("I116").Value = 2841.8
$ws.Range("J116").Value = 9502.833000000001
$ws.Range("K116").Value = 2841.8
$ws.Range("L116").Value = 9502.833000000001
$ws.Range("M116").Value = -547.8000000000002
$ws.Range("N116").Value = -14090.833

$ws.Range("H122").Value = 2979.9443
$ws.Range("I122").Value = 2735.5386
$ws.Range("K122").Value = 8206.6158
$ws.Range("M122").Value = -5756.6158

$ws.Range("H132").Value = 1186
$ws.Range("I132").Value = 1186
$ws.Range("K132").Value = 3558
$ws.Range("M132").Value = -1028

$ws.Range("H136").Value = 906.3333
$ws.Range("I136").Value = 609.5
$ws.Range("K136").Value = 1828.5
$ws.Range("M136").Value = 721.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6475.091
$ws.Range("I3").Value = 2841.8
$ws.Range("J3").Value = 9502.833000000001
$ws.Range("K3").Value = 2841.8
$ws.Range("L3").Value = 9502.833000000001
$ws.Range("M3").Value = -2727.8
$ws.Range("N3").Value = -9730.833000000001

$ws.Range("H99").Value = 4934.7144
$ws.Range("I99").Value = 3610.75
$ws.Range("K99").Value = 3610.75
$ws.Range("M99").Value = -2112.75

$ws.Range("H105").Value = 4218.0625
$ws.Range("I105").Value = 3999.2307
$ws.Range("K105").Value = 3999.2307
$ws.Range("M105").Value = -2252.2307

$ws.Range("H134").Value = 1187.8
$ws.Range("I134").Value = 986.44446
$ws.Range("K134").Value = 2959.33338
$ws.Range("M134").Value = -424.33338

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2543.25
$ws.Range("I31").Value = 2637
$ws.Range("J31").Value = 2449.5
$ws.Range("K31").Value = 2637
$ws.Range("L31").Value = 2449.5
$ws.Range("M31").Value = -2342
$ws.Range("N31").Value = -3039.5

$ws.Range("H34").Value = 2543.25
$ws.Range("I34").Value = 2637
$ws.Range("J34").Value = 2449.5
$ws.Range("K34").Value = 2637
$ws.Range("L34").Value = 2449.5
$ws.Range("M34").Value = -2435
$ws.Range("N34").Value = -2853.5

$ws.Range("H105").Value = 2030.2
$ws.Range("I105").Value = 2060.6
$ws.Range("J105").Value = 1999.8
$ws.Range("K105").Value = 2060.6
$ws.Range("L105").Value = 1999.8
$ws.Range("M105").Value = -313.5999999999999
$ws.Range("N105").Value = -5493.8

$ws.Range("H122").Value = 848
$ws.Range("I122").Value = 774.2857
$ws.Range("J122").Value = 934
$ws.Range("K122").Value = 2322.8571
$ws.Range("L122").Value = 2802
$ws.Range("M122").Value = 127.1428999999998
$ws.Range("N122").Value = -7702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1234
$ws.Range("I5").Value = 1234
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3702
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -3590

$ws.Range("H135").Value = 1234
$ws.Range("I135").Value = 1234
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 11106
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2333.5715
$ws.Range("I80").Value = 2199.6667
$ws.Range("K80").Value = 2199.6667
$ws.Range("M80").Value = -1201.6667

$ws.Range("H83").Value = 2333.5715
$ws.Range("I83").Value = 2199.6667
$ws.Range("K83").Value = 10998.3335
$ws.Range("M83").Value = -6006.333500000001

$ws.Range("H102").Value = 1191.75
$ws.Range("I102").Value = 1191.75
$ws.Range("K102").Value = 1191.75
$ws.Range("M102").Value = 430.25

$ws.Range("H122").Value = 1616.6666
$ws.Range("I122").Value = 2233.3333
$ws.Range("K122").Value = 6699.999899999999
$ws.Range("M122").Value = -4249.999899999999

$ws.Range("H123").Value = 27000
$ws.Range("J123").Value = 27000
$ws.Range("L123").Value = 27000
$ws.Range("N123").Value = -31900

$ws.Range("H125").Value = 65999.5
$ws.Range("J125").Value = 65999.5
$ws.Range("L125").Value = 65999.5
$ws.Range("N125").Value = -70919.5

$ws.Range("H126").Value = 4222.25
$ws.Range("J126").Value = 4599
$ws.Range("L126").Value = 13797
$ws.Range("N126").Value = -18737

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705

$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893

$ws.Range("H46").Value = 1589.5
$ws.Range("I46").Value = 1250
$ws.Range("J46").Value = 1929
$ws.Range("K46").Value = 1250
$ws.Range("L46").Value = 1929
$ws.Range("M46").Value = -1062
$ws.Range("N46").Value = -2305

$ws.Range("H55").Value = 271.92
$ws.Range("I55").Value = 174.4
$ws.Range("J55").Value = 296.3
$ws.Range("K55").Value = 174.4
$ws.Range("L55").Value = 296.3
$ws.Range("M55").Value = -1.400000000000006
$ws.Range("N55").Value = -642.3

$ws.Range("H68").Value = 2038.5555
$ws.Range("I68").Value = 2058.3333
$ws.Range("J68").Value = 1999
$ws.Range("K68").Value = 2058.3333
$ws.Range("L68").Value = 1999
$ws.Range("M68").Value = -1309.3333
$ws.Range("N68").Value = -3497

$ws.Range("H71").Value = 2038.5555
$ws.Range("I71").Value = 2058.3333
$ws.Range("J71").Value = 1999
$ws.Range("K71").Value = 10291.6665
$ws.Range("L71").Value = 9995
$ws.Range("M71").Value = -6547.666499999999
$ws.Range("N71").Value = -17483

$ws.Range("H122").Value = 4881.846
$ws.Range("I122").Value = 4797
$ws.Range("K122").Value = 14391
$ws.Range("M122").Value = -11941

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1617
$ws.Range("I96").Value = 1269.5714
$ws.Range("K96").Value = 1269.5714
$ws.Range("M96").Value = 103.4286

$ws.Range("H122").Value = 1181.8096
$ws.Range("J122").Value = 1166.4445
$ws.Range("L122").Value = 3499.3335
$ws.Range("N122").Value = -8399.333500000001

$ws.Range("H136").Value = 1801.1666
$ws.Range("I136").Value = 1368
$ws.Range("J136").Value = 2523.111
$ws.Range("K136").Value = 4104
$ws.Range("L136").Value = 7569.333
$ws.Range("M136").Value = -1554
$ws.Range("N136").Value = -12669.333
